# Updated cryptos list on Sun Jan  7 04:46:40 UTC 2024 with GitHub Actions
# Refreshes the Price (col D) / Volume(1h) (col E) snapshot for each coin row,
# and re-ranks a few coins whose order changed (B/C swap together with D/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.441.39'
$ws.Range("E2").Value = '  +0.59%  '

# Row 3
$ws.Range("D3").Value = '2.248.54'
$ws.Range("E3").Value = '  -0.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").Value = '307.56'
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
$ws.Range("D6").Value = '94.82'
$ws.Range("E6").Value = '  -3.71%  '

# Row 7
$ws.Range("E7").Value = '  -0.91%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.23%  '

# Row 9
$ws.Range("E9").Value = '  -1.19%  '

# Row 10
$ws.Range("D10").Value = '34.84'
$ws.Range("E10").Value = '  -2.08%  '

# Row 11
$ws.Range("E11").Value = '  -1.00%  '

# Row 12
$ws.Range("D12").Value = '7.21'
$ws.Range("E12").Value = '  -1.34%  '

# Row 13
$ws.Range("E13").Value = '  +0.62%  '

# Row 14
$ws.Range("D14").Value = '2.411.00'
$ws.Range("E14").Value = '  +5.26%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.840'
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = '  -0.55%  '

# Row 17
$ws.Range("D17").Value = '44.104.42'
$ws.Range("E17").Value = '  +0.21%  '

# Row 18
$ws.Range("D18").Value = '12.51'
$ws.Range("E18").Value = '  -2.36%  '

# Row 19
$ws.Range("E19").Value = '  -1.12%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.38%  '

# Row 21
$ws.Range("E21").Value = '  +0.97%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '237.59'
$ws.Range("E22").Value = '  -1.74%  '

# Row 23
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '3.01'
$ws.Range("E23").Value = '  +2.04%  '

# Row 24
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  +2.18%  '

# Row 25
$ws.Range("E25").Value = '  -0.21%  '

# Row 26
$ws.Range("D26").Value = '38.38'
$ws.Range("E26").Value = '  +4.82%  '

# Row 27
$ws.Range("E27").Value = '  +3.71%  '

# Row 28
$ws.Range("D28").Value = '9.89'
$ws.Range("E28").Value = '  -2.14%  '

# Row 29
$ws.Range("E29").Value = '  -3.64%  '

# Row 30
$ws.Range("D30").Value = '20.12'
$ws.Range("E30").Value = '  -0.14%  '

# Row 31
$ws.Range("D31").Value = '154.02'
$ws.Range("E31").Value = '  -1.83%  '

# Row 32
$ws.Range("D32").Value = '0.0803'
$ws.Range("E32").Value = '  -2.12%  '

# Row 33
$ws.Range("E33").Value = '  -0.59%  '

# Row 34
$ws.Range("E34").Value = '  -11.32%  '

# Row 35
$ws.Range("D35").Value = '0.109'
$ws.Range("E35").Value = '  +2.03%  '

# Row 36
$ws.Range("E36").Value = '  +0.43%  '

# Row 37
$ws.Range("D37").Value = '1.81'
$ws.Range("E37").Value = '  -1.77%  '

# Row 38
$ws.Range("E38").Value = '  +3.59%  '

# Row 39
$ws.Range("D39").Value = '14.84'
$ws.Range("E39").Value = '  -4.06%  '

# Row 40
$ws.Range("D40").Value = '3.83'
$ws.Range("E40").Value = '  -0.80%  '

# Row 41
$ws.Range("D41").Value = '0.0304'
$ws.Range("E41").Value = '  -0.74%  '

# Row 42
$ws.Range("E42").Value = '  +0.34%  '

# Row 43
$ws.Range("D43").Value = '1.749.96'
$ws.Range("E43").Value = '  -0.64%  '

# Row 44
$ws.Range("E44").Value = '  +0.94%  '

# Row 45
$ws.Range("D45").Value = '81.02'
$ws.Range("E45").Value = '  -6.70%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '99.85'
$ws.Range("E46").Value = '  -1.44%  '

# Row 47
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").Value = '4.96'
$ws.Range("E47").Value = '  -3.90%  '

# Row 48
$ws.Range("D48").Value = '70.83'
$ws.Range("E48").Value = '  +1.93%  '

# Row 49
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '1.62'
$ws.Range("E49").Value = '  +5.69%  '

# Row 50
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '56.04'
$ws.Range("E50").Value = '  +0.87%  '

# Row 51
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '8.16'
$ws.Range("E51").Value = '  -1.23%  '
